# Insert a new row at row 55 (this pushes former rows 55-168 down to 56-169)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with its data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the same values as the row that used to
# be at position 55 (now at 56); only D, J, K, L, M, P differ.
$ws.Cells.Item(55, 1).Value  = 4
$ws.Cells.Item(55, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(55, 3).Value  = "Los Lagos"
$ws.Cells.Item(55, 4).Value  = 44469
$ws.Cells.Item(55, 5).Value  = 10
$ws.Cells.Item(55, 6).Value  = 100112040
$ws.Cells.Item(55, 7).Value  = "Cilantro"
$ws.Cells.Item(55, 8).Value  = "Sin especificar"
$ws.Cells.Item(55, 9).Value  = "Primera"
$ws.Cells.Item(55, 10).Value = 150
$ws.Cells.Item(55, 11).Value = 11000
$ws.Cells.Item(55, 12).Value = 11000
$ws.Cells.Item(55, 13).Value = 11000
$ws.Cells.Item(55, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(55, 15).Value = "Región Metropolitana"
$ws.Cells.Item(55, 16).Value = 306
$ws.Cells.Item(55, 17).Value = 36
$ws.Cells.Item(55, 18).Value = "Hortaliza"
